$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-117 down to 46-118
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the new record
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C45").Value = 'Arica y Parinacota'
$ws.Range("D45").Value = 44792
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 'Fruta'
$ws.Range("G45").Value = 100106
$ws.Range("H45").Value = 'Oleaginosos'
$ws.Range("I45").Value = 100106002
$ws.Range("J45").Value = 'Palta'
$ws.Range("K45").Value = 'Hass'
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 400
$ws.Range("N45").Value = 23000
$ws.Range("O45").Value = 24000
$ws.Range("P45").Value = 23500
$ws.Range("Q45").Value = '$/bandeja 10 kilos'
$ws.Range("R45").Value = 'Perú'
$ws.Range("S45").Value = 2350
$ws.Range("T45").Value = 10
